# Updated fish survey datasheet for IC-C1 transects 1 and 2 (2020.12.14).
# Adds the newly reviewed frames (rows 43-83) for transect 2, backfilling
# the date/site/transect/frame/fileName columns on the rows that already
# existed (43-45) and appending entirely new rows (46-83).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

  # Row 43
  $ws.Cells.Item(43, 1).Value = 44179
  $ws.Cells.Item(43, 2).Value = "IC-C1"
  $ws.Cells.Item(43, 3).Value = 2
  $ws.Cells.Item(43, 4).Value = 3
  $ws.Cells.Item(43, 5).Value = "GOPR0195.MP4007.png"
  $ws.Cells.Item(43, 6).Value = "French Grunt"
  $ws.Cells.Item(43, 7).Value = "Haemulon"
  $ws.Cells.Item(43, 8).Value = "flavolineatum"
  $ws.Cells.Item(43, 9).Value = "Haemulon flavolineatum"
  $ws.Cells.Item(43, 10).Value = 5
  # Row 44
  $ws.Cells.Item(44, 1).Value = 44179
  $ws.Cells.Item(44, 2).Value = "IC-C1"
  $ws.Cells.Item(44, 3).Value = 2
  $ws.Cells.Item(44, 4).Value = 3
  $ws.Cells.Item(44, 5).Value = "GOPR0195.MP4007.png"
  $ws.Cells.Item(44, 6).Value = "Caesar Grunt"
  $ws.Cells.Item(44, 7).Value = "Haemulon"
  $ws.Cells.Item(44, 8).Value = "carbonarium"
  $ws.Cells.Item(44, 9).Value = "Haemulon carbonarium"
  $ws.Cells.Item(44, 10).Value = 1
  # Row 45
  $ws.Cells.Item(45, 1).Value = 44179
  $ws.Cells.Item(45, 2).Value = "IC-C1"
  $ws.Cells.Item(45, 3).Value = 2
  $ws.Cells.Item(45, 4).Value = 3
  $ws.Cells.Item(45, 5).Value = "GOPR0195.MP4007.png"
  $ws.Cells.Item(45, 6).Value = "Dusky Damselfish"
  $ws.Cells.Item(45, 7).Value = "Stegastes"
  $ws.Cells.Item(45, 8).Value = "adustus"
  $ws.Cells.Item(45, 9).Value = "Stegastes adustus"
  $ws.Cells.Item(45, 10).Value = 2
  # Row 46
  $ws.Cells.Item(46, 1).Value = 44179
  $ws.Cells.Item(46, 2).Value = "IC-C1"
  $ws.Cells.Item(46, 3).Value = 2
  $ws.Cells.Item(46, 4).Value = 3
  $ws.Cells.Item(46, 5).Value = "GOPR0195.MP4007.png"
  $ws.Cells.Item(46, 6).Value = "Bicolor Damselfish"
  $ws.Cells.Item(46, 7).Value = "Stegastes"
  $ws.Cells.Item(46, 8).Value = "partitus"
  $ws.Cells.Item(46, 9).Value = "Stegastes partitus"
  $ws.Cells.Item(46, 10).Value = 1
  # Row 47
  $ws.Cells.Item(47, 1).Value = 44179
  $ws.Cells.Item(47, 2).Value = "IC-C1"
  $ws.Cells.Item(47, 3).Value = 2
  $ws.Cells.Item(47, 4).Value = 4
  $ws.Cells.Item(47, 5).Value = "GOPR0195.MP4008.png"
  $ws.Cells.Item(47, 6).Value = "Bicolor Damselfish"
  $ws.Cells.Item(47, 7).Value = "Stegastes"
  $ws.Cells.Item(47, 8).Value = "partitus"
  $ws.Cells.Item(47, 9).Value = "Stegastes partitus"
  $ws.Cells.Item(47, 10).Value = 1
  # Row 48
  $ws.Cells.Item(48, 1).Value = 44179
  $ws.Cells.Item(48, 2).Value = "IC-C1"
  $ws.Cells.Item(48, 3).Value = 2
  $ws.Cells.Item(48, 4).Value = 4
  $ws.Cells.Item(48, 5).Value = "GOPR0195.MP4008.png"
  $ws.Cells.Item(48, 6).Value = "French Grunt"
  $ws.Cells.Item(48, 7).Value = "Haemulon "
  $ws.Cells.Item(48, 8).Value = "flavolineatum"
  $ws.Cells.Item(48, 9).Value = "Haemulon flavolineatum"
  $ws.Cells.Item(48, 10).Value = 1
  # Row 49
  $ws.Cells.Item(49, 1).Value = 44179
  $ws.Cells.Item(49, 2).Value = "IC-C1"
  $ws.Cells.Item(49, 3).Value = 2
  $ws.Cells.Item(49, 4).Value = 4
  $ws.Cells.Item(49, 5).Value = "GOPR0195.MP4008.png"
  $ws.Cells.Item(49, 6).Value = "Bluehead"
  $ws.Cells.Item(49, 7).Value = "Thalassoma"
  $ws.Cells.Item(49, 8).Value = "bifasciatum"
  $ws.Cells.Item(49, 9).Value = "Thalassoma bifasciatum"
  $ws.Cells.Item(49, 10).Value = 1
  $ws.Cells.Item(49, 12).Value = "juvenile"
  $ws.Cells.Item(49, 13).Value = "yellow morph"
  # Row 50
  $ws.Cells.Item(50, 1).Value = 44179
  $ws.Cells.Item(50, 2).Value = "IC-C1"
  $ws.Cells.Item(50, 3).Value = 2
  $ws.Cells.Item(50, 4).Value = 5
  $ws.Cells.Item(50, 5).Value = "GOPR0195.MP4009.png"
  $ws.Cells.Item(50, 6).Value = "Caesar Grunt"
  $ws.Cells.Item(50, 7).Value = "Haemulon"
  $ws.Cells.Item(50, 8).Value = "carbonarium"
  $ws.Cells.Item(50, 9).Value = "Haelumon carbonarium"
  $ws.Cells.Item(50, 10).Value = 1
  # Row 51
  $ws.Cells.Item(51, 1).Value = 44179
  $ws.Cells.Item(51, 2).Value = "IC-C1"
  $ws.Cells.Item(51, 3).Value = 2
  $ws.Cells.Item(51, 4).Value = 6
  $ws.Cells.Item(51, 5).Value = "GOPR0195.MP4010.png"
  $ws.Cells.Item(51, 6).Value = "Bicolor Damselfish"
  $ws.Cells.Item(51, 7).Value = "Stegastes"
  $ws.Cells.Item(51, 8).Value = "partitus"
  $ws.Cells.Item(51, 9).Value = "Stegastes partitus"
  $ws.Cells.Item(51, 10).Value = 2
  # Row 52
  $ws.Cells.Item(52, 1).Value = 44179
  $ws.Cells.Item(52, 2).Value = "IC-C1"
  $ws.Cells.Item(52, 3).Value = 2
  $ws.Cells.Item(52, 4).Value = 6
  $ws.Cells.Item(52, 5).Value = "GOPR0195.MP4010.png"
  $ws.Cells.Item(52, 6).Value = "Bluehead"
  $ws.Cells.Item(52, 7).Value = "Thalassoma"
  $ws.Cells.Item(52, 8).Value = "bifasciatum"
  $ws.Cells.Item(52, 9).Value = "Thalassoma bifasciatum"
  $ws.Cells.Item(52, 10).Value = 2
  $ws.Cells.Item(52, 12).Value = "juvenile"
  $ws.Cells.Item(52, 13).Value = "yellow morph"
  # Row 53
  $ws.Cells.Item(53, 1).Value = 44179
  $ws.Cells.Item(53, 2).Value = "IC-C1"
  $ws.Cells.Item(53, 3).Value = 2
  $ws.Cells.Item(53, 4).Value = 7
  $ws.Cells.Item(53, 5).Value = "GOPR0195.MP4011.png"
  $ws.Cells.Item(53, 6).Value = "Bicolor Damselfish"
  $ws.Cells.Item(53, 7).Value = "Stegastes "
  $ws.Cells.Item(53, 8).Value = "partitus"
  $ws.Cells.Item(53, 9).Value = "Stegastes partitus"
  $ws.Cells.Item(53, 10).Value = 4
  # Row 54
  $ws.Cells.Item(54, 1).Value = 44179
  $ws.Cells.Item(54, 2).Value = "IC-C1"
  $ws.Cells.Item(54, 3).Value = 2
  $ws.Cells.Item(54, 4).Value = 7
  $ws.Cells.Item(54, 5).Value = "GOPR0195.MP4011.png"
  $ws.Cells.Item(54, 6).Value = "Spanish Hogfish"
  $ws.Cells.Item(54, 7).Value = "Bodianus"
  $ws.Cells.Item(54, 8).Value = "rufus"
  $ws.Cells.Item(54, 9).Value = "Bodianus rufus"
  $ws.Cells.Item(54, 10).Value = 1
  # Row 55
  $ws.Cells.Item(55, 1).Value = 44179
  $ws.Cells.Item(55, 2).Value = "IC-C1"
  $ws.Cells.Item(55, 3).Value = 2
  $ws.Cells.Item(55, 4).Value = 7
  $ws.Cells.Item(55, 5).Value = "GOPR0195.MP4011.png"
  $ws.Cells.Item(55, 6).Value = "Bluehead"
  $ws.Cells.Item(55, 7).Value = "Thalassoma"
  $ws.Cells.Item(55, 8).Value = "bifasciatum"
  $ws.Cells.Item(55, 9).Value = "Thalassoma bifasciatum"
  $ws.Cells.Item(55, 10).Value = 2
  # Row 56
  $ws.Cells.Item(56, 1).Value = 44179
  $ws.Cells.Item(56, 2).Value = "IC-C1"
  $ws.Cells.Item(56, 3).Value = 2
  $ws.Cells.Item(56, 4).Value = 7
  $ws.Cells.Item(56, 5).Value = "GOPR0195.MP4011.png"
  $ws.Cells.Item(56, 6).Value = "Foureye Butterflyfish"
  $ws.Cells.Item(56, 7).Value = "Cheatodon"
  $ws.Cells.Item(56, 8).Value = "capistratus"
  $ws.Cells.Item(56, 9).Value = "Chaetodon capistratus"
  $ws.Cells.Item(56, 10).Value = 1
  # Row 57
  $ws.Cells.Item(57, 1).Value = 44179
  $ws.Cells.Item(57, 2).Value = "IC-C1"
  $ws.Cells.Item(57, 3).Value = 2
  $ws.Cells.Item(57, 4).Value = 7
  $ws.Cells.Item(57, 5).Value = "GOPR0195.MP4011.png"
  $ws.Cells.Item(57, 6).Value = "Painted Wrasse"
  $ws.Cells.Item(57, 7).Value = "Halichoeres"
  $ws.Cells.Item(57, 8).Value = "caudalis"
  $ws.Cells.Item(57, 9).Value = "Halichoeres caudalis"
  $ws.Cells.Item(57, 10).Value = 1
  # Row 58
  $ws.Cells.Item(58, 1).Value = 44179
  $ws.Cells.Item(58, 2).Value = "IC-C1"
  $ws.Cells.Item(58, 3).Value = 2
  $ws.Cells.Item(58, 4).Value = 8
  $ws.Cells.Item(58, 5).Value = "GOPR0195.MP4012.png"
  $ws.Cells.Item(58, 6).Value = "Painted Wrasse"
  $ws.Cells.Item(58, 7).Value = "Halichoeres"
  $ws.Cells.Item(58, 8).Value = "caudalis"
  $ws.Cells.Item(58, 9).Value = "Halichoeres caudalis"
  $ws.Cells.Item(58, 10).Value = 1
  # Row 59
  $ws.Cells.Item(59, 1).Value = 44179
  $ws.Cells.Item(59, 2).Value = "IC-C1"
  $ws.Cells.Item(59, 3).Value = 2
  $ws.Cells.Item(59, 4).Value = 8
  $ws.Cells.Item(59, 5).Value = "GOPR0195.MP4012.png"
  $ws.Cells.Item(59, 6).Value = "Foureye butterflyfish"
  $ws.Cells.Item(59, 7).Value = "Chaetodon"
  $ws.Cells.Item(59, 8).Value = "capistratus"
  $ws.Cells.Item(59, 9).Value = "Chaetodon capistratus"
  $ws.Cells.Item(59, 10).Value = 1
  # Row 60
  $ws.Cells.Item(60, 1).Value = 44179
  $ws.Cells.Item(60, 2).Value = "IC-C1"
  $ws.Cells.Item(60, 3).Value = 2
  $ws.Cells.Item(60, 4).Value = 8
  $ws.Cells.Item(60, 5).Value = "GOPR0195.MP4012.png"
  $ws.Cells.Item(60, 6).Value = "Bicolor Damselfish"
  $ws.Cells.Item(60, 7).Value = "Stegastes"
  $ws.Cells.Item(60, 8).Value = "partitus"
  $ws.Cells.Item(60, 9).Value = "Stegastes partitus"
  $ws.Cells.Item(60, 10).Value = 2
  # Row 61
  $ws.Cells.Item(61, 1).Value = 44179
  $ws.Cells.Item(61, 2).Value = "IC-C1"
  $ws.Cells.Item(61, 3).Value = 2
  $ws.Cells.Item(61, 4).Value = 8
  $ws.Cells.Item(61, 5).Value = "GOPR0195.MP4012.png"
  $ws.Cells.Item(61, 6).Value = "Bluehead"
  $ws.Cells.Item(61, 7).Value = "Thalassoma"
  $ws.Cells.Item(61, 8).Value = "bifasciatum"
  $ws.Cells.Item(61, 9).Value = "Thalassoma bifasciatum"
  $ws.Cells.Item(61, 10).Value = 2
  $ws.Cells.Item(61, 12).Value = "juvenile"
  $ws.Cells.Item(61, 13).Value = "yellow morph"
  # Row 62
  $ws.Cells.Item(62, 1).Value = 44179
  $ws.Cells.Item(62, 2).Value = "IC-C1"
  $ws.Cells.Item(62, 3).Value = 2
  $ws.Cells.Item(62, 4).Value = 9
  $ws.Cells.Item(62, 5).Value = "GOPR0195.MP4013.png"
  $ws.Cells.Item(62, 6).Value = "Bluehead"
  $ws.Cells.Item(62, 7).Value = "Thalassoma"
  $ws.Cells.Item(62, 8).Value = "bifasciatum"
  $ws.Cells.Item(62, 9).Value = "Thalassoma bifasciatum"
  $ws.Cells.Item(62, 10).Value = 2
  $ws.Cells.Item(62, 12).Value = "juvenile"
  $ws.Cells.Item(62, 13).Value = "yellow morph"
  # Row 63
  $ws.Cells.Item(63, 1).Value = 44179
  $ws.Cells.Item(63, 2).Value = "IC-C1"
  $ws.Cells.Item(63, 3).Value = 2
  $ws.Cells.Item(63, 4).Value = 9
  $ws.Cells.Item(63, 5).Value = "GOPR0195.MP4013.png"
  $ws.Cells.Item(63, 6).Value = "Painted Wrasse"
  $ws.Cells.Item(63, 7).Value = "Halichoeres"
  $ws.Cells.Item(63, 8).Value = "caudalis"
  $ws.Cells.Item(63, 9).Value = "Halichoeres caudalis"
  $ws.Cells.Item(63, 10).Value = 1
  # Row 64
  $ws.Cells.Item(64, 1).Value = 44179
  $ws.Cells.Item(64, 2).Value = "IC-C1"
  $ws.Cells.Item(64, 3).Value = 2
  $ws.Cells.Item(64, 4).Value = 10
  $ws.Cells.Item(64, 5).Value = "GOPR0195.MP4014.png"
  $ws.Cells.Item(64, 6).Value = "N/A"
  $ws.Cells.Item(64, 7).Value = "N/A"
  $ws.Cells.Item(64, 8).Value = "N/A"
  $ws.Cells.Item(64, 9).Value = "N/A"
  $ws.Cells.Item(64, 10).Value = "N/A"
  # Row 65
  $ws.Cells.Item(65, 1).Value = 44179
  $ws.Cells.Item(65, 2).Value = "IC-C1"
  $ws.Cells.Item(65, 3).Value = 2
  $ws.Cells.Item(65, 4).Value = 11
  $ws.Cells.Item(65, 5).Value = "GOPR0195.MP4015.png"
  $ws.Cells.Item(65, 6).Value = "Bicolor Damselfish"
  $ws.Cells.Item(65, 7).Value = "Stegastes"
  $ws.Cells.Item(65, 8).Value = "partitus"
  $ws.Cells.Item(65, 9).Value = "Stegastes partitus"
  $ws.Cells.Item(65, 10).Value = 7
  # Row 66
  $ws.Cells.Item(66, 1).Value = 44179
  $ws.Cells.Item(66, 2).Value = "IC-C1"
  $ws.Cells.Item(66, 3).Value = 2
  $ws.Cells.Item(66, 4).Value = 11
  $ws.Cells.Item(66, 5).Value = "GOPR0195.MP4015.png"
  $ws.Cells.Item(66, 6).Value = "Painted Wrasse"
  $ws.Cells.Item(66, 7).Value = "Halichoeres"
  $ws.Cells.Item(66, 8).Value = "caudalis"
  $ws.Cells.Item(66, 9).Value = "Halichoeres caudalis"
  $ws.Cells.Item(66, 10).Value = 1
  # Row 67
  $ws.Cells.Item(67, 1).Value = 44179
  $ws.Cells.Item(67, 2).Value = "IC-C1"
  $ws.Cells.Item(67, 3).Value = 2
  $ws.Cells.Item(67, 4).Value = 11
  $ws.Cells.Item(67, 5).Value = "GOPR0195.MP4015.png"
  $ws.Cells.Item(67, 6).Value = "Bluehead"
  $ws.Cells.Item(67, 7).Value = "Thalassoma"
  $ws.Cells.Item(67, 8).Value = "bifasciatum"
  $ws.Cells.Item(67, 9).Value = "Thalassoma bifasciatum"
  $ws.Cells.Item(67, 10).Value = 5
  # Row 68
  $ws.Cells.Item(68, 1).Value = 44179
  $ws.Cells.Item(68, 2).Value = "IC-C1"
  $ws.Cells.Item(68, 3).Value = 2
  $ws.Cells.Item(68, 4).Value = 12
  $ws.Cells.Item(68, 5).Value = "GOPR0195.MP4016.png"
  $ws.Cells.Item(68, 6).Value = "Bicolor Damselfish"
  $ws.Cells.Item(68, 7).Value = "Stegastes"
  $ws.Cells.Item(68, 8).Value = "partitus"
  $ws.Cells.Item(68, 9).Value = "Stegastes partitus"
  $ws.Cells.Item(68, 10).Value = 5
  # Row 69
  $ws.Cells.Item(69, 1).Value = 44179
  $ws.Cells.Item(69, 2).Value = "IC-C1"
  $ws.Cells.Item(69, 3).Value = 2
  $ws.Cells.Item(69, 4).Value = 12
  $ws.Cells.Item(69, 5).Value = "GOPR0195.MP4016.png"
  $ws.Cells.Item(69, 6).Value = "Bluehead"
  $ws.Cells.Item(69, 7).Value = "Thalassoma"
  $ws.Cells.Item(69, 8).Value = "bifasciatum"
  $ws.Cells.Item(69, 9).Value = "Thalassoma bifasciatum"
  $ws.Cells.Item(69, 10).Value = 9
  $ws.Cells.Item(69, 12).Value = "juvenile"
  $ws.Cells.Item(69, 13).Value = "yellow morph"
  # Row 70
  $ws.Cells.Item(70, 1).Value = 44179
  $ws.Cells.Item(70, 2).Value = "IC-C1"
  $ws.Cells.Item(70, 3).Value = 2
  $ws.Cells.Item(70, 4).Value = 13
  $ws.Cells.Item(70, 5).Value = "GOPR0195.MP4017.png"
  $ws.Cells.Item(70, 6).Value = "Bluehead"
  $ws.Cells.Item(70, 7).Value = "Thalassoma"
  $ws.Cells.Item(70, 8).Value = "bifasciatum"
  $ws.Cells.Item(70, 9).Value = "Thalassoma bifasciatum"
  $ws.Cells.Item(70, 10).Value = 12
  $ws.Cells.Item(70, 12).Value = "juvenile"
  $ws.Cells.Item(70, 13).Value = "yellow morph"
  # Row 71
  $ws.Cells.Item(71, 1).Value = 44179
  $ws.Cells.Item(71, 2).Value = "IC-C1"
  $ws.Cells.Item(71, 3).Value = 2
  $ws.Cells.Item(71, 4).Value = 14
  $ws.Cells.Item(71, 5).Value = "GOPR0195.MP4018.png"
  $ws.Cells.Item(71, 6).Value = "Bluestriped Grunt"
  $ws.Cells.Item(71, 7).Value = "Haemulon"
  $ws.Cells.Item(71, 8).Value = "sciurus"
  $ws.Cells.Item(71, 9).Value = "Haemulon sciurus"
  $ws.Cells.Item(71, 10).Value = 1
  # Row 72
  $ws.Cells.Item(72, 1).Value = 44179
  $ws.Cells.Item(72, 2).Value = "IC-C1"
  $ws.Cells.Item(72, 3).Value = 2
  $ws.Cells.Item(72, 4).Value = 14
  $ws.Cells.Item(72, 5).Value = "GOPR0195.MP4018.png"
  $ws.Cells.Item(72, 6).Value = "Bicolor Damselfish"
  $ws.Cells.Item(72, 7).Value = "Stegastes"
  $ws.Cells.Item(72, 8).Value = "partitus"
  $ws.Cells.Item(72, 9).Value = "Segastes partitus"
  $ws.Cells.Item(72, 10).Value = 5
  # Row 73
  $ws.Cells.Item(73, 1).Value = 44179
  $ws.Cells.Item(73, 2).Value = "IC-C1"
  $ws.Cells.Item(73, 3).Value = 2
  $ws.Cells.Item(73, 4).Value = 14
  $ws.Cells.Item(73, 5).Value = "GOPR0195.MP4018.png"
  $ws.Cells.Item(73, 6).Value = "Bluestriped Grunt"
  $ws.Cells.Item(73, 7).Value = "Haemulon"
  $ws.Cells.Item(73, 8).Value = "sciurus"
  $ws.Cells.Item(73, 9).Value = "Haemulon sciurus"
  $ws.Cells.Item(73, 10).Value = 1
  # Row 74
  $ws.Cells.Item(74, 1).Value = 44179
  $ws.Cells.Item(74, 2).Value = "IC-C1"
  $ws.Cells.Item(74, 3).Value = 2
  $ws.Cells.Item(74, 4).Value = 14
  $ws.Cells.Item(74, 5).Value = "GOPR0195.MP4018.png"
  $ws.Cells.Item(74, 6).Value = "Bluehead"
  $ws.Cells.Item(74, 7).Value = "Thalassoma"
  $ws.Cells.Item(74, 8).Value = "bifasciatum"
  $ws.Cells.Item(74, 9).Value = "Thalassoma bifasciatum"
  $ws.Cells.Item(74, 10).Value = 6
  # Row 75
  $ws.Cells.Item(75, 1).Value = 44179
  $ws.Cells.Item(75, 2).Value = "IC-C1"
  $ws.Cells.Item(75, 3).Value = 2
  $ws.Cells.Item(75, 4).Value = 14
  $ws.Cells.Item(75, 5).Value = "GOPR0195.MP4018.png"
  $ws.Cells.Item(75, 6).Value = "Slippery Dick "
  $ws.Cells.Item(75, 7).Value = "Halichoeres "
  $ws.Cells.Item(75, 8).Value = "bivittatus"
  $ws.Cells.Item(75, 9).Value = "Halichoeres bivittatus"
  $ws.Cells.Item(75, 10).Value = 2
  $ws.Cells.Item(75, 12).Value = "juvenile"
  # Row 76
  $ws.Cells.Item(76, 1).Value = 44179
  $ws.Cells.Item(76, 2).Value = "IC-C1"
  $ws.Cells.Item(76, 3).Value = 2
  $ws.Cells.Item(76, 4).Value = 15
  $ws.Cells.Item(76, 5).Value = "GOPR0195.MP4019.png"
  $ws.Cells.Item(76, 6).Value = "Bluehead"
  $ws.Cells.Item(76, 7).Value = "Thalassoma"
  $ws.Cells.Item(76, 8).Value = "bifasciatum"
  $ws.Cells.Item(76, 9).Value = "Thalassoma bifasciatum"
  $ws.Cells.Item(76, 10).Value = 4
  # Row 77
  $ws.Cells.Item(77, 1).Value = 44179
  $ws.Cells.Item(77, 2).Value = "IC-C1"
  $ws.Cells.Item(77, 3).Value = 2
  $ws.Cells.Item(77, 4).Value = 16
  $ws.Cells.Item(77, 5).Value = "GOPR0195.MP4020.png"
  $ws.Cells.Item(77, 6).Value = "Bicolor Damselfish"
  $ws.Cells.Item(77, 7).Value = "Stegastes "
  $ws.Cells.Item(77, 8).Value = "partitus"
  $ws.Cells.Item(77, 9).Value = "Stegastes partitus"
  $ws.Cells.Item(77, 10).Value = 3
  # Row 78
  $ws.Cells.Item(78, 1).Value = 44179
  $ws.Cells.Item(78, 2).Value = "IC-C1"
  $ws.Cells.Item(78, 3).Value = 2
  $ws.Cells.Item(78, 4).Value = 16
  $ws.Cells.Item(78, 5).Value = "GOPR0195.MP4020.png"
  $ws.Cells.Item(78, 6).Value = "White Grunt"
  $ws.Cells.Item(78, 7).Value = "Haemulon"
  $ws.Cells.Item(78, 8).Value = "plumierii"
  $ws.Cells.Item(78, 9).Value = "Haemulon pulmierii"
  $ws.Cells.Item(78, 10).Value = 1
  # Row 79
  $ws.Cells.Item(79, 1).Value = 44179
  $ws.Cells.Item(79, 2).Value = "IC-C1"
  $ws.Cells.Item(79, 3).Value = 2
  $ws.Cells.Item(79, 4).Value = 16
  $ws.Cells.Item(79, 5).Value = "GOPR0195.MP4020.png"
  $ws.Cells.Item(79, 6).Value = "Parrotfish"
  $ws.Cells.Item(79, 7).Value = "UNK"
  $ws.Cells.Item(79, 8).Value = "UNK"
  $ws.Cells.Item(79, 9).Value = "UNK"
  $ws.Cells.Item(79, 10).Value = 2
  $ws.Cells.Item(79, 13).Value = "both individuals unidentifiable "
  # Row 80
  $ws.Cells.Item(80, 1).Value = 44179
  $ws.Cells.Item(80, 2).Value = "IC-C1"
  $ws.Cells.Item(80, 3).Value = 2
  $ws.Cells.Item(80, 4).Value = 17
  $ws.Cells.Item(80, 5).Value = "GOPR0195.MP4021.png"
  $ws.Cells.Item(80, 6).Value = "White Grunt"
  $ws.Cells.Item(80, 7).Value = "Haemulon"
  $ws.Cells.Item(80, 8).Value = "plumierii"
  $ws.Cells.Item(80, 9).Value = "Haemulon pulmierii"
  $ws.Cells.Item(80, 10).Value = 1
  # Row 81
  $ws.Cells.Item(81, 1).Value = 44179
  $ws.Cells.Item(81, 2).Value = "IC-C1"
  $ws.Cells.Item(81, 3).Value = 2
  $ws.Cells.Item(81, 4).Value = 17
  $ws.Cells.Item(81, 5).Value = "GOPR0195.MP4021.png"
  $ws.Cells.Item(81, 6).Value = "Parrotfish"
  $ws.Cells.Item(81, 7).Value = "UNK"
  $ws.Cells.Item(81, 8).Value = "UNK"
  $ws.Cells.Item(81, 9).Value = "UNK"
  $ws.Cells.Item(81, 10).Value = 1
  $ws.Cells.Item(81, 13).Value = "individual unidentifiable"
  # Row 82
  $ws.Cells.Item(82, 1).Value = 44179
  $ws.Cells.Item(82, 2).Value = "IC-C1"
  $ws.Cells.Item(82, 3).Value = 2
  $ws.Cells.Item(82, 4).Value = 18
  $ws.Cells.Item(82, 5).Value = "GOPR0195.MP4022.png"
  $ws.Cells.Item(82, 6).Value = "White Grunt"
  $ws.Cells.Item(82, 7).Value = "Haemulon"
  $ws.Cells.Item(82, 8).Value = "plumierii"
  $ws.Cells.Item(82, 9).Value = "Haemulon pulmierii"
  $ws.Cells.Item(82, 10).Value = 1
  # Row 83
  $ws.Cells.Item(83, 1).Value = 44179
  $ws.Cells.Item(83, 2).Value = "IC-C1"
  $ws.Cells.Item(83, 3).Value = 2
  $ws.Cells.Item(83, 4).Value = 18
  $ws.Cells.Item(83, 5).Value = "GOPR0195.MP4022.png"
  $ws.Cells.Item(83, 6).Value = "Parrotfish"
  $ws.Cells.Item(83, 7).Value = "UNK"
  $ws.Cells.Item(83, 8).Value = "UNK"
  $ws.Cells.Item(83, 9).Value = "UNK"
  $ws.Cells.Item(83, 10).Value = 1
  $ws.Cells.Item(83, 13).Value = "individual unidentifiable "

# Give the new date cells (column A) the same date number format as the
# rest of the sheet (style index carried over from an existing date cell,
# so no new number format gets minted).
$ws.Range("A2").Copy()
$ws.Range("A43:A83").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Restore the view: scrolled/selected state the author left the sheet in
# after finishing the 2020.12.14 entry pass.
$ws.Activate()
$ws.Range("M83").Select()

Write-Output "Updated rows 43-83 on $($ws.Name)."
